$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.257.09"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "3.541.74"
$ws.Range("E3").Value = "  -3.51%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'610.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.51%  "
$ws.Range("D6").Value = "'154.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D7").Value = "3.539.28"
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.484"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").Value = "'6.86"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").Value = "'0.429"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'0.0000223"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").Value = "4.141.61"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").Value = "'31.98"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "3.540.28"
$ws.Range("E16").Value = "  -3.06%  "
$ws.Range("D17").Value = "67.198.92"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'6.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'15.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "'446.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("D22").Value = "'9.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.78%  "
$ws.Range("D23").Value = "'0.632"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'78.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "3.679.02"
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'0.0000123"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "'10.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("D29").Value = "'8.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.17%  "
$ws.Range("D30").Value = "'2.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").Value = "'1.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'25.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("D35").Value = "'1.88"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'6.19"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.533.90"
$ws.Range("E37").Value = "  -3.51%  "
$ws.Range("D38").Value = "'8.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'174.99"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "'2.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").Value = "'5.58"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.54%  "
$ws.Range("D44").Value = "'0.0868"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").Value = "'0.893"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").Value = "'45.74"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "'27.64"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'2.63"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'1.23"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").Value = "'7.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "'1.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.51%  "
